# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# These two sheets mirror the same event data, so both need identical updates.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 8
    3 = 1769
    4 = 547
    5 = 1124
    6 = 5974
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Range("F$row").Value = $newValues[$row]
    }
}
